$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

# Copy existing header style (bold + border + center/top alignment) onto A1
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New "hole_id" index column header
$ws.Range("A1").Value = "hole_id"

# hole_id values for rows 2-47, replacing the old numeric 0..45 index
$holeIds = @(
    "LBU_05_02",
    "LBU_07_02",
    "LBU_05_07",
    "MHZ_08_01",
    "MHZ_08_05",
    "MHZ_12_04",
    "LBU_96_4",
    "LBU_05_29",
    "LBU_05_14",
    "LBU_05_13",
    "LBU_05_26",
    "LBU_05_09",
    "LBU_96_3",
    "LBU_05_05",
    "MHZ_12_03",
    "LBU_98_4",
    "LBU_07_01",
    "LBU_87_5",
    "MHZ_12_01",
    "LBU_05_03",
    "LBU_87_2",
    "LBU_02_3",
    "LBU_07_03",
    "LBU_02_4",
    "LBU_05_20",
    "LBU_98_5",
    "LBU_05_15",
    "LBU_05_24",
    "LBU_05_12",
    "LBU_87_4",
    "LBU_05_27",
    "LBU_05_18",
    "LBU_05_30",
    "LBU_87_3",
    "MHZ_08_04",
    "MHZ_08_02",
    "LBU_87_6",
    "LBU_96_2",
    "LBU_05_08",
    "LBU_05_17",
    "LBU_05_19",
    "LBU_01_1",
    "LBU_05_06",
    "LBU_05_01",
    "LBU_05_22",
    "LBU_96_1"
)

for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $holeIds[$i]
}

